# Applies the OOXML changes described by the diff:
#  - M1: update the "80..." regex note text (swap text previously
#    held by the shared string now used by C1 / M1)
#  - M2:M22: fill in the "x" marker cells that were missing in column M
#  - Update the active selection from D24 to L25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note text in M1 (the "80 ..." regex description), keeping
# the same leading issue number and trailing sentence but refreshing
# the regex line in the middle.
$ws.Range("M1").Value2 = "80`n\b(?!(base|utils|grDevices|graphics)\b)\w+(?=::)`ndo it after colons_check"

# Fill the previously-empty M column (rows 2-22) with the same "x"
# marker used throughout the rest of the table.
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 13).Value2 = "x"
}

# Move the active selection to L25 (was D24).
$ws.Range("L25").Select()
